# Active_Outages.xlsx update - 6/19/2025, 5:52:44 AM
# Bumps elapsed-duration strings for each open outage and converts the
# orphaned "MAK0875" entry on sheet R5 into a proper follow-up row.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3943:06:11"
$ws1.Range("G3").Value = "82:38:49"
$ws1.Range("G4").Value = "105:38:49"

$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12124:29:52"
$ws2.Range("G3").Value = "3254:13:21"
$ws2.Range("G4").Value = "492:24:55"

$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2970:19:41"
$ws4.Range("G3").Value = "197:31:56"
$ws4.Range("G4").Value = "85:44:21"
$ws4.Range("G5").Value = "83:21:54"

$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "444:19:08"

# Row 3 used to just hold a leftover "MAK0875" hub-site value in column A.
# Replace it with a full follow-up row describing the still-open MAK0605 outage.
$ws5.Range("A3").Value = ""
$ws5.Range("B3").Value = "R5"
$ws5.Range("C3").Value = ""
$ws5.Range("D3").Value = "MAK0605"
$ws5.Range("E3").Value = ""
$ws5.Range("F3").Value = ""
$ws5.Range("G3").Value = ""
$ws5.Range("H3").Value = ""
$ws5.Range("I3").Value = "SCECO"
$ws5.Range("J3").Value = "Dead"
$ws5.Range("K3").Value = ""
$ws5.Range("L3").Value = "Latis"

$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "84:50:58"
